# "Generate Report for Archive"
# - Update the localization status text from "Ready for handoff" to
#   "In Translation" everywhere it appears (Overview!E2/F2, zh-cn!C2,
#   de-de!C2).
# - Narrow the status columns that held that text (Overview columns E & F,
#   and column C on the zh-cn / de-de sheets) to match the shorter text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# --- Update the status text -------------------------------------------------
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$zhcn.Range("C2").Value = $newStatus
$dede.Range("C2").Value = $newStatus

# --- Narrow the affected status columns -------------------------------------
# ColumnWidth of 12.5 characters is the closest settable width to the
# target stored column width of ~13.41 in this engine's column-width grid.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
